$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 1 values (A1:L1 integers, M1:Q1 fractional "percentage-like" values) ----
$ws.Range("A1").Value = 3
$ws.Range("B1").Value = 2
$ws.Range("C1").Value = 18
$ws.Range("D1").Value = 12
$ws.Range("E1").Value = 16
$ws.Range("F1").Value = 6
$ws.Range("G1").Value = 20
$ws.Range("H1").Value = 13
$ws.Range("I1").Value = 31
$ws.Range("J1").Value = 6
$ws.Range("K1").Value = 19
$ws.Range("L1").Value = 16
$ws.Range("M1").Value = 1.7/20
$ws.Range("N1").Value = 0.037
$ws.Range("O1").Value = 0.031
$ws.Range("P1").Value = 0.073
$ws.Range("Q1").Value = 0.7/20

# ---- Column widths (columns A-E are unchanged from the source file) ----
$ws.Columns.Item(6).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(7).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(8).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(9).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(10).ColumnWidth = 1.3333333333333333
$ws.Columns.Item(11).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(12).ColumnWidth = 2.3333333333333335
$ws.Columns.Item(13).ColumnWidth = 4.833333333333333
$ws.Columns.Item(14).ColumnWidth = 4.833333333333333
$ws.Columns.Item(15).ColumnWidth = 4.833333333333333
$ws.Columns.Item(16).ColumnWidth = 4.833333333333333
$ws.Columns.Item(17).ColumnWidth = 4.833333333333333
